$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'36.708.64"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").Value = "'1.957.16"
$ws.Range("E3").Value = "  +1.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'244.74"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +2.21%  "

# Row 7
$ws.Range("D7").Value = "'59.01"
$ws.Range("E7").Value = "  +5.85%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.370"
$ws.Range("E9").Value = "  +2.91%  "

# Row 10
$ws.Range("D10").Value = "'0.0814"
$ws.Range("E10").Value = "  -0.21%  "

# Row 11
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("D12").Value = "'22.68"
$ws.Range("E12").Value = "  +9.81%  "

# Row 13
$ws.Range("D13").Value = "'2.239.72"
$ws.Range("E13").Value = "  +1.66%  "

# Row 14
$ws.Range("D14").Value = "'0.823"
$ws.Range("E14").Value = "  +2.38%  "

# Row 15
$ws.Range("D15").Value = "'13.71"
$ws.Range("E15").Value = "  +4.62%  "

# Row 16
$ws.Range("D16").Value = "'5.27"
$ws.Range("E16").Value = "  +2.58%  "

# Row 17
$ws.Range("D17").Value = "'1.965.68"
$ws.Range("E17").Value = "  +2.02%  "

# Row 18
$ws.Range("D18").Value = "'36.665.41"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19
$ws.Range("D19").Value = "'69.84"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20
$ws.Range("E20").Value = "  +1.20%  "

# Row 21
$ws.Range("D21").Value = "'228.55"
$ws.Range("E21").Value = "  +1.74%  "

# Row 22
$ws.Range("D22").Value = "'5.07"
$ws.Range("E22").Value = "  +3.23%  "

# Row 23
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("D24").Value = "'2.42"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  +3.59%  "

# Row 26
$ws.Range("D26").Value = "'9.33"
$ws.Range("E26").Value = "  +1.23%  "

# Row 27
$ws.Range("D27").Value = "'160.56"
$ws.Range("E27").Value = "  -0.88%  "

# Row 28
$ws.Range("D28").Value = "'0.136"
$ws.Range("E28").Value = "  +17.43%  "

# Row 29
$ws.Range("D29").Value = "'19.37"
$ws.Range("E29").Value = "  +1.80%  "

# Row 30
$ws.Range("E30").Value = "  +2.64%  "

# Row 31
$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32
$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = "  +2.51%  "

# Row 33
$ws.Range("D33").Value = "'0.0623"
$ws.Range("E33").Value = "  +1.60%  "

# Row 34
$ws.Range("D34").Value = "'4.24"
$ws.Range("E34").Value = "  +1.00%  "

# Row 35
$ws.Range("D35").Value = "'6.26"
$ws.Range("E35").Value = "  +5.78%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'3.48"
$ws.Range("E36").Value = "  +22.39%  "

# Row 37
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38
$ws.Range("D38").Value = "'2.22"
$ws.Range("E38").Value = "  +4.81%  "

# Row 39
$ws.Range("D39").Value = "'1.76"
$ws.Range("E39").Value = "  -0.95%  "

# Row 40
$ws.Range("E40").Value = "  +5.55%  "

# Row 41
$ws.Range("E41").Value = "  +2.65%  "

# Row 42
$ws.Range("D42").Value = "'0.0212"
$ws.Range("E42").Value = "  +3.54%  "

# Row 43
$ws.Range("D43").Value = "'1.16"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44
$ws.Range("D44").Value = "'16.12"
$ws.Range("E44").Value = "  +4.89%  "

# Row 45
$ws.Range("E45").Value = "  +2.59%  "

# Row 46
$ws.Range("D46").Value = "'1.349.71"
$ws.Range("E46").Value = "  +1.65%  "

# Row 47
$ws.Range("D47").Value = "'87.63"
$ws.Range("E47").Value = "  +1.40%  "

# Row 48
$ws.Range("D48").Value = "'7.21"
$ws.Range("E48").Value = "  +1.27%  "

# Row 49
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  +1.47%  "

# Row 50
$ws.Range("D50").Value = "'2.131.66"
$ws.Range("E50").Value = "  +1.71%  "

# Row 51
$ws.Range("D51").Value = "'43.57"
$ws.Range("E51").Value = "  -3.44%  "
